$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("P3").Value = "Condition.code"
$ws.Range("S3").Value = ""
$ws.Range("P4").Value = ".stage.summary"
$ws.Range("S4").Value = "Is required."
$ws.Range("P5").Value = ".onsetDateTime"
$ws.Range("S5").Value = ""
$ws.Range("P8").Value = ".dateOfLastDressingChange(extension, datetime)"
$ws.Range("P9").Value = "Observation(pressureUlcurObservable).component.woundLength (valueQuantity, ucum)"
$ws.Range("P10").Value = "Observation(pressureUlcurObservable).component.woundWidth (valueQuantity, ucum)"
$ws.Range("P11").Value = "Observation(pressureUlcurObservable).component.woundDepth (valueQuantity, ucum)"
$ws.Range("P12").Value = "Media (resource)"
$ws.Range("P13").Value = ".note"
